# Append new data row (row 70) to the sheet reflecting the 2026-02-02 profit run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/02/2026"
$ws.Cells.Item($row, 2).Value = 9441.6
$ws.Cells.Item($row, 3).Value = 0.2646812310596689
$ws.Cells.Item($row, 4).Value = 0.7353187689403311
$ws.Cells.Item($row, 5).Value = -329.55
$ws.Cells.Item($row, 6).Value = -41.98
$ws.Cells.Item($row, 7).Value = -23899.18
$ws.Cells.Item($row, 8).Value = -77.48999999999999
$ws.Cells.Item($row, 9).Value = -803.84
$ws.Cells.Item($row, 10).Value = -24.34
$ws.Cells.Item($row, 11).Value = -24703.02
$ws.Cells.Item($row, 12).Value = -72.34999999999999
